# LAB_LinkedArt_ModelTracking.xlsx
# - Remove <xsl:text> usage note -> add a "replaced line breaks" note in F12
# - Reformat the Dimension Statement row (17): mark mapped, fill URI/Type/Notes
# - Add dimension-statement URI and AAT vocab term
# - Update sheet view: drop stale topLeftCell on the view, re-home the frozen
#   pane to A2, and move the live selection to B14

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objects")

# Row 17 (Dimension Statement): mark as mapped and complete the mapping info
$ws.Range("B17").Value = "X"
$ws.Range("D17").Value = "object/[irn]/dimension-statement"
$ws.Range("E17").Value = "http://vocab.getty.edu/aat/300266036"

# Row 12 (Title Notes): add the clean-up note about replacing line breaks
$ws.Range("F12").Value = "Replaced line breaks with replace(., '\n', '\\n')"
$ws.Range("F17").Value = "Replaced line breaks with replace(., '\n', '\\n')"

# Refresh the view: re-anchor the frozen pane/selection
[void]$ws.Range("B14").Select()
